$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Forandrad" (column C) date for every existing data row (2-121)
#    from 2023-10-05 (45204) to 2023-10-06 (45205).
$ws.Range("C2:C121").Value2 = 45205

# 2. Row 121 gains an explicit row height (15) / customHeight flag.
$ws.Rows.Item(121).RowHeight = 15

# 3. Copy the per-cell formatting used on row 121 (date format for B/C, wrap-text
#    style for R) down onto the new row 122, cell by cell, so we don't drag along
#    any formatting for the otherwise-empty columns (e.g. F) like a full-row copy would.
$ws.Cells.Item(121, 2).Copy($ws.Cells.Item(122, 2))
$ws.Cells.Item(121, 3).Copy($ws.Cells.Item(122, 3))
$ws.Cells.Item(121, 18).Copy($ws.Cells.Item(122, 18))
$excel.CutCopyMode = 0

# 4. Fill in the new entry "A 47860-2023" on row 122.
$ws.Cells.Item(122, 1).Value = "A 47860-2023"
$ws.Cells.Item(122, 2).Value2 = 45204
$ws.Cells.Item(122, 3).Value2 = 45205
$ws.Cells.Item(122, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(122, 5).Value = "MULLSJÖ"
$ws.Cells.Item(122, 7).Value2 = 1.6
$ws.Cells.Item(122, 8).Value2 = 0
$ws.Cells.Item(122, 9).Value2 = 0
$ws.Cells.Item(122, 10).Value2 = 0
$ws.Cells.Item(122, 11).Value2 = 0
$ws.Cells.Item(122, 12).Value2 = 0
$ws.Cells.Item(122, 13).Value2 = 0
$ws.Cells.Item(122, 14).Value2 = 0
$ws.Cells.Item(122, 15).Value2 = 0
$ws.Cells.Item(122, 16).Value2 = 0
$ws.Cells.Item(122, 17).Value2 = 0
